{"js": "// The document contains four \"<id>...</id>\" tags, each originally split\n// across three runs: \"<id>\", the bare id text (e.g. \"p007r_a1\"), and\n// \"</id>\". The edit collapses each of those three runs into a single run\n// whose text is \"<id>p007r_N</id>\" (the \"a\" prefix is dropped from the\n// numeric suffix), while keeping the surrounding document untouched.\nconst oldIds = [\"p007r_a1\", \"p007r_a2\", \"p007r_a3\", \"p007r_a4\"];\nconst newIds = [\"p007r_1\", \"p007r_2\", \"p007r_3\", \"p007r_4\"];\n\nfor (let i = 0; i < oldIds.length; i++) {\n  const searchResults = context.document.body.search(`<id>${oldIds[i]}</id>`, {\n    matchCase: true,\n  });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (let j = 0; j < searchResults.items.length; j++) {\n    searchResults.items[j].insertText(`<id>${newIds[i]}</id>`, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains four \"<id>...</id>\" tags, each originally split\n# across three runs: \"<id>\", the bare id text (e.g. \"p007r_a1\"), and\n# \"</id>\". The edit collapses each of those three runs into a single run\n# whose text is \"<id>p007r_N</id>\" (the \"a\" prefix is dropped from the\n# numeric suffix), leaving the rest of the document untouched.\n$d = $word.ActiveDocument\n\n$oldIds = @(\"p007r_a1\", \"p007r_a2\", \"p007r_a3\", \"p007r_a4\")\n$newIds = @(\"p007r_1\", \"p007r_2\", \"p007r_3\", \"p007r_4\")\n\nfor ($i = 0; $i -lt $oldIds.Length; $i++) {\n    $searchText = \"<id>$($oldIds[$i])</id>\"\n    $replaceText = \"<id>$($newIds[$i])</id>\"\n\n    $range = $d.Content\n    $range.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
